# Rename the "Window Number" header (A1) to "Window Type".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Window Type"

# Reset the sheet's selection back to the top-left cell (the saved file
# had a stray selection at O9 left over from editing).
$ws.Range("A1").Select() | Out-Null
